# Re-run of the "100000" grid-search cell (Output_14_14.xlsx): updates the
# cached results in Summary, Fed-in Capacity, Unmet Demand, Household
# Surplus and Costs and Revenues to the newly solved values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = 898.4058559126861
$ws.Range("B7").Value = 10414372.90082427
$ws.Range("B8").Value = 24690005.52504999
$ws.Range("B10").Value = 3057988.557334842

$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("M14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("Q20").Value = 150.3839754851235
$ws.Range("J21").Value = 93.17061249236157
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = 113.4004983079896
$ws.Range("O23").Value = 117.8828208804077
$ws.Range("I24").Value = 10.12574714858493
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("J26").Value = 124.5190384721106
$ws.Range("L26").Value = 130.6648563030561
$ws.Range("O26").Value = 117.8828208804077
$ws.Range("P26").Value = 135.4597561231036
$ws.Range("Q26").Value = 150.3839754851235
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 51.84373129681028
$ws.Range("O27").Value = 57.81213424001893
$ws.Range("P27").Value = 0
$ws.Range("L28").Value = 90.4687457914608
$ws.Range("N28").Value = 81.96869489115805
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("Q32").Value = 150.3839754851235
$ws.Range("J33").Value = 93.17061249236157
$ws.Range("Q33").Value = 0
$ws.Range("Q35").Value = 150.3839754851235
$ws.Range("Q36").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("J38").Value = 124.5190384721106
$ws.Range("J40").Value = 33.63624132272333
$ws.Range("L40").Value = 90.4687457914608
$ws.Range("M40").Value = 92.09541281912071
$ws.Range("J41").Value = 124.5190384721106
$ws.Range("Q42").Value = 94.49434172313325
$ws.Range("N43").Value = 81.96869489115805

$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("M14").Value = 113.4004983079896
$ws.Range("O14").Value = 117.8828208804077
$ws.Range("K15").Value = 80.29914934735042
$ws.Range("L15").Value = 61.18167021676314
$ws.Range("M15").Value = 51.84373129681028
$ws.Range("N15").Value = 38.66169381481656
$ws.Range("O15").Value = 57.81213424001893
$ws.Range("Q15").Value = 94.49434172313325
$ws.Range("L16").Value = 90.4687457914608
$ws.Range("M16").Value = 92.09541281912071
$ws.Range("N16").Value = 81.96869489115805
$ws.Range("O16").Value = 96.22962838366004
$ws.Range("K17").Value = 135.370731907559
$ws.Range("L17").Value = 130.6648563030561
$ws.Range("M17").Value = 113.4004983079896
$ws.Range("O17").Value = 117.8828208804077
$ws.Range("P17").Value = 135.4597561231036
$ws.Range("J18").Value = 93.17061249236157
$ws.Range("L18").Value = 61.18167021676314
$ws.Range("N18").Value = 38.66169381481656
$ws.Range("O18").Value = 57.81213424001893
$ws.Range("P18").Value = 65.92768427608706
$ws.Range("K19").Value = 94.30397654773019
$ws.Range("L19").Value = 90.4687457914608
$ws.Range("M19").Value = 92.09541281912071
$ws.Range("N19").Value = 81.96869489115805
$ws.Range("K20").Value = 135.370731907559
$ws.Range("N20").Value = 110.5750244233121
$ws.Range("Q20").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("N22").Value = 81.96869489115805
$ws.Range("O22").Value = 96.22962838366004
$ws.Range("K23").Value = 135.370731907559
$ws.Range("M23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("I24").Value = 77.12765456497084
$ws.Range("O24").Value = 57.81213424001893
$ws.Range("P24").Value = 65.92768427608706
$ws.Range("Q24").Value = 94.49434172313325
$ws.Range("R24").Value = 123.5547069419379
$ws.Range("K25").Value = 94.30397654773019
$ws.Range("L25").Value = 90.4687457914608
$ws.Range("M25").Value = 92.09541281912071
$ws.Range("O25").Value = 96.22962838366004
$ws.Range("P25").Value = 101.5955875616828
$ws.Range("Q25").Value = 126.4887893424616
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("K27").Value = 80.29914934735042
$ws.Range("L27").Value = 61.18167021676314
$ws.Range("M27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 65.92768427608706
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("K30").Value = 80.29914934735042
$ws.Range("L30").Value = 61.18167021676314
$ws.Range("P30").Value = 65.92768427608706
$ws.Range("N32").Value = 110.5750244233121
$ws.Range("Q32").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("Q33").Value = 94.49434172313325
$ws.Range("Q35").Value = 0
$ws.Range("Q36").Value = 94.49434172313325
$ws.Range("L37").Value = 90.4687457914608
$ws.Range("N37").Value = 81.96869489115805
$ws.Range("J38").Value = 0
$ws.Range("J40").Value = 72.23757736389061
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("N43").Value = 0

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B6").Value = 288809.8323486014
$ws.Range("B7").Value = 234126.9916772344
$ws.Range("B8").Value = 204358.8141148659
$ws.Range("B9").Value = 230454.468325944
$ws.Range("B10").Value = 315235.943839665
$ws.Range("B11").Value = 234750.1002025046
$ws.Range("B12").Value = 246595.763001371
$ws.Range("B13").Value = 220496.2475379781
$ws.Range("B14").Value = 213801.6124174268
$ws.Range("B15").Value = 193864.4916710733

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("F2").Value = 64847.52161019725
$ws.Range("G2").Value = 53910.95347592385
$ws.Range("H2").Value = 47957.31796345015
$ws.Range("I2").Value = 53176.44880566579
$ws.Range("J2").Value = 70132.74390840999
$ws.Range("K2").Value = 54035.57518097788
$ws.Range("L2").Value = 56404.70774075118
$ws.Range("M2").Value = 51184.80464807258
$ws.Range("N2").Value = 49845.87762396233
$ws.Range("O2").Value = 45858.45347469164
$ws.Range("E3").Value = 133100.0000000001
$ws.Range("F4").Value = 29419.74590781255
$ws.Range("G4").Value = 18483.17777353916
$ws.Range("H4").Value = 12529.54226106547
$ws.Range("I4").Value = 17748.67310328111
$ws.Range("J4").Value = 34704.9682060253
$ws.Range("K4").Value = 18607.7994785932
$ws.Range("L4").Value = 20976.93203836649
$ws.Range("M4").Value = 15757.0289456879
$ws.Range("N4").Value = 14418.10192157765
$ws.Range("O4").Value = 10430.67777230695
$ws.Range("E6").Value = -101283.8364026237
$ws.Range("F6").Value = 31816.16359737637
$ws.Range("H6").Value = 31816.16359737636
$ws.Range("J6").Value = 31816.16359737636
